# The "Recorded By" column (G) lists the people/processes who recorded
# attendance for a session, as a comma-separated string. Previously the
# literal "System" entry (when present) was listed first; it should
# instead be listed last, e.g. "System, user@x.com" -> "user@x.com, System".
# Only entries where "System" (capital S) is the first of multiple
# comma-separated values are affected - single-value "System" cells and
# cells where "System" already isn't first are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$prefix = "System, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text
    if ($v.StartsWith($prefix)) {
        $rest = $v.Substring($prefix.Length)
        $cell.Value = $rest + ", System"
    }
}
